$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

$ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3)).Value = 45205
